$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.527.64"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "2.478.46"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.20"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.12"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  -1.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +1.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.66"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0784"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("E12").Value = "  +2.41%  "

$ws.Range("D13").Value = "2.862.17"
$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("E14").Value = "  -1.88%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.15"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +8.65%  "

$ws.Range("D16").Value = "2.456.51"
$ws.Range("E16").Value = "  +0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.764"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.05%  "

$ws.Range("D18").Value = "41.507.74"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("E19").Value = "  +2.26%  "

$ws.Range("D20").Value = "0.0₃0935"
$ws.Range("E20").Value = "  +2.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.51"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.27"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.74%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.77"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.70"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.80%  "

$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.68%  "

$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.65"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.80%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.05"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.86"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.46"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0753"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.04%  "

$ws.Range("B35").Value = "ApeXProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.44"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -7.60%  "

$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.31"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.91"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.84%  "

$ws.Range("E38").Value = "  +3.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.82"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.14"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.73%  "

$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.43"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.06%  "

$ws.Range("D44").Value = "1.967.12"
$ws.Range("E44").Value = "  -0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0283"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.95"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.02"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.95%  "

$ws.Range("D48").Value = "2.719.59"
$ws.Range("E48").Value = "  +0.89%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.71"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.93"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.87%  "
